$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices for B,C,D,E,F,I,J,K,L,M,N (H is intentionally skipped/blank)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

# New vm_pu values for rows 2..25 (one row per data row), columns in $cols order
$data = @(
    @(1.02,1.03593698750209,1.042376069805158,1.048926541866687,1.055187971965981,1.034200394410536,1.041048038162894,1.045152932583596,1.051685025127399,1.057929117875372,1.017610318966543),
    @(1.02,1.036853001531433,1.043063750560743,1.049796540051587,1.056066075677844,1.034321600963353,1.04160797694474,1.045651755013737,1.052367017186028,1.058620448208254,1.017798480631754),
    @(1.02,1.037446408948675,1.04350911122926,1.050360492172367,1.056635213038827,1.034398796860583,1.04197034835639,1.045974233307898,1.052808675002547,1.059068089554908,1.017920186398779),
    @(1.02,1.037696040117263,1.043696431417422,1.050597816603435,1.056874703015495,1.034430954255274,1.042122700950344,1.046109731749914,1.052994433540726,1.059256349378941,1.017971339700079),
    @(1.02,1.037737963790832,1.043727888563547,1.050637678383905,1.056914927558384,1.034436336261961,1.042148282270443,1.046132478309454,1.053025628195289,1.059287963176081,1.017979927859151),
    @(1.02,1.037449743894952,1.043511613857024,1.050363662376176,1.056638412235151,1.034399227712168,1.041972384054893,1.04597604412696,1.052811156782658,1.059070604813176,1.017920869958641),
    @(1.02,1.036246416084429,1.042608393971097,1.049220352706584,1.055484534402783,1.034241611632089,1.041237260296055,1.045321571913363,1.051915431463577,1.058162692287238,1.017673918627285),
    @(1.02,1.034131304176424,1.041019829043407,1.047213466771506,1.053458572650583,1.03395445699935,1.039942340552499,1.044166117980776,1.050339893507493,1.056565230709326,1.017238415696462),
    @(1.02,1.032724874821993,1.039962917772947,1.045880867282077,1.052112955100405,1.033756727040704,1.039079441203681,1.043394416806253,1.049291530207883,1.055501954837132,1.016947876568971),
    @(1.02,1.03211675625135,1.039505790028176,1.045305120711513,1.051531501049963,1.033669622705145,1.038705901107579,1.043059945507665,1.04883806797594,1.055041965854423,1.016822027207195),
    @(1.02,1.031891006621841,1.039336072500767,1.045091456449397,1.051315706240619,1.033637045624505,1.03856716802285,1.042935661197267,1.0486697063174,1.054871169504819,1.016775274989254),
    @(1.02,1.031939424658821,1.039372473836282,1.045137279368998,1.05136198664548,1.033644043584777,1.038596925982427,1.042962322700652,1.04870581707538,1.054907802980101,1.016785303770316),
    @(1.02,1.032098093002402,1.039491759481323,1.045287455196884,1.051513659638156,1.033666934412078,1.038694433047485,1.043049673074889,1.048824149622851,1.055027846460009,1.016818162777371),
    @(1.02,1.032195871425974,1.039565265964348,1.045380009191367,1.051607134691305,1.033681008720568,1.038754512568244,1.043103486361912,1.048897068111751,1.055101817739678,1.016838407493076),
    @(1.02,1.032765252185323,1.039993266972766,1.045919104750485,1.052151569874343,1.033762476612023,1.039104234057119,1.043416607936916,1.049321635360722,1.055532491691094,1.016956227887915),
    @(1.02,1.033122644757253,1.040261881583689,1.046257608437438,1.052493404119664,1.033813181855312,1.039323632972692,1.043612936137026,1.049588086274417,1.055802754623384,1.017030122045395),
    @(1.02,1.033331190099823,1.04041861017651,1.046455175088089,1.052692906614787,1.033842613952336,1.039451614316806,1.043727420134961,1.04974354938261,1.055960434505777,1.017073218981555),
    @(1.02,1.033402312932486,1.04047205904617,1.046522561052434,1.052760951494253,1.03385262519328,1.039495254252431,1.043766450937564,1.049796566180507,1.056014206025595,1.017087913184155),
    @(1.02,1.033084291187594,1.040233056569958,1.046221277442521,1.052456716491656,1.033807756492205,1.039300092536954,1.043591875168398,1.049559493762445,1.055773753811561,1.017022194340635),
    @(1.02,1.032051365456869,1.039456630608194,1.04524322679949,1.051468990653323,1.03366019977547,1.038665719186754,1.043023951848178,1.048789301571881,1.05499249486759,1.016808486785264),
    @(1.02,1.031402692361077,1.038968924881415,1.044629408765854,1.050849028521141,1.033566137110871,1.038266958590832,1.04266660588724,1.048305482344607,1.054501657601492,1.016674084936745),
    @(1.02,1.031746492986253,1.039227422365611,1.044954698481035,1.051177581091052,1.033616123387717,1.038478339724067,1.042856066924464,1.048561922736824,1.054761824043843,1.016745337107634),
    @(1.02,1.033101621254458,1.040246081197886,1.046237693474872,1.052473293692211,1.033810208423451,1.039310729416831,1.043601391809237,1.049572413343412,1.055786857906404,1.017025776542435),
    @(1.02,1.034677474544414,1.041430142608781,1.047731363679726,1.0539814544314,1.034029805314178,1.040277046389037,1.044465082446056,1.050746862212095,1.056977919843576,1.017351041537499)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = 2 + $r
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($rowNum, $cols[$c]).Value = $rowVals[$c]
    }
}

Write-Host "Done updating vm_pu values"
